$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: dpc / substance_id row (sus_id) ---
# B4 text content is unchanged ("substance_id"); shared-string reindex only, no action needed.

# --- Row 25: dpc_id -----------------------------------------------------
# Un-highlight (remove the red "Bad" style) and fill in the new mapping.
$ws.Range("A25:D25").ClearFormats()
$ws.Range("B25").Value = "coordinate_precision_id"
$ws.Range("C25").Value = """coordinate_precisions"""
# D25 keeps its text ("Precision of coordinates"); only formatting changed.

# --- Row 32: dic_id -------------------------------------------------------
$ws.Range("B32").Value = "concentration_indicator_id"
# D32 text unchanged ("concentration indicator").

# --- Rows 38-41: sampling_date_y/m/d/t ------------------------------------
$ws.Range("B38").Value = "sampling_date_year"
$ws.Range("B39").Value = "sampling_date_month"
$ws.Range("B40").Value = "sampling_date_date"
$ws.Range("B41").Value = "sampling_date_time"

# --- Rows 47-48: remark / remark_add --------------------------------------
$ws.Range("C47").Value = "moved to ""remarks"" json"
$ws.Range("C48").Value = "moved to ""remarks"" json"

# --- Row 58: dtl_id --------------------------------------------------------
$ws.Range("A58:D58").ClearFormats()
$ws.Range("B58").Value = "treatment_less_id"
$ws.Range("C58").Value = "treatment_less"
$ws.Range("D58").ClearContents()

# --- Row 59: dtl_other -------------------------------------------------
$ws.Range("A59:D59").ClearFormats()
$ws.Range("B59").ClearContents()
$ws.Range("C59").Value = "moved to ""treatment_less"""
$ws.Range("D59").ClearContents()

# --- Row 60: sampling_date1 (B60 dropped) ---------------------------------
$ws.Range("B60").ClearContents()

# --- Row 61: sampling_date1_y (B61 dropped) -------------------------------
$ws.Range("B61").ClearContents()

# --- Row 71: list_id -------------------------------------------------------
$ws.Range("A71:D71").ClearFormats()
$ws.Range("B71").Value = "file_source_id"
$ws.Range("C71").ClearContents()
$ws.Range("D71").ClearContents()

# --- Selection / view state -------------------------------------------------
$ws.Range("A14").Select()
$excel.ActiveWindow.Zoom = 145
